# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Reorders the worker/period data table (rows 16-43) from being grouped by
# worker (periods 2311 -> 2305 descending) to being grouped by period
# (2305 -> 2311 ascending, with the same 4 workers repeated for each
# period). Also updates the "Valor Mora" (column G) from 1160000 to
# 1000000 for every row, and the "Salario Basico" (column F) so the new
# last period (2311) carries the 37333 figure that used to belong to the
# first row of each worker's block.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$rows = @(
    @{Row=16; C="60405011";   D="ROSARIO MARIA FLOREZ TEHERAN";   E="2305"; F=46400; G=1000000},
    @{Row=17; C="1047455587"; D="FERNANDO ANDRES OVALLE CORDERO"; E="2305"; F=46400; G=1000000},
    @{Row=18; C="7920858";    D="RUBEN DARIO CORONEL MORALES";    E="2305"; F=46400; G=1000000},
    @{Row=19; C="1065003738"; D="JUAN MANUEL DEGIOVANNI PRECIADO";E="2305"; F=46400; G=1000000},
    @{Row=20; C="60405011";   D="ROSARIO MARIA FLOREZ TEHERAN";   E="2306"; F=46400; G=1000000},
    @{Row=21; C="1047455587"; D="FERNANDO ANDRES OVALLE CORDERO"; E="2306"; F=46400; G=1000000},
    @{Row=22; C="7920858";    D="RUBEN DARIO CORONEL MORALES";    E="2306"; F=46400; G=1000000},
    @{Row=23; C="1065003738"; D="JUAN MANUEL DEGIOVANNI PRECIADO";E="2306"; F=46400; G=1000000},
    @{Row=24; C="60405011";   D="ROSARIO MARIA FLOREZ TEHERAN";   E="2307"; F=46400; G=1000000},
    @{Row=25; C="1047455587"; D="FERNANDO ANDRES OVALLE CORDERO"; E="2307"; F=46400; G=1000000},
    @{Row=26; C="7920858";    D="RUBEN DARIO CORONEL MORALES";    E="2307"; F=46400; G=1000000},
    @{Row=27; C="1065003738"; D="JUAN MANUEL DEGIOVANNI PRECIADO";E="2307"; F=46400; G=1000000},
    @{Row=28; C="60405011";   D="ROSARIO MARIA FLOREZ TEHERAN";   E="2308"; F=46400; G=1000000},
    @{Row=29; C="1047455587"; D="FERNANDO ANDRES OVALLE CORDERO"; E="2308"; F=46400; G=1000000},
    @{Row=30; C="7920858";    D="RUBEN DARIO CORONEL MORALES";    E="2308"; F=46400; G=1000000},
    @{Row=31; C="1065003738"; D="JUAN MANUEL DEGIOVANNI PRECIADO";E="2308"; F=46400; G=1000000},
    @{Row=32; C="60405011";   D="ROSARIO MARIA FLOREZ TEHERAN";   E="2309"; F=46400; G=1000000},
    @{Row=33; C="1047455587"; D="FERNANDO ANDRES OVALLE CORDERO"; E="2309"; F=46400; G=1000000},
    @{Row=34; C="7920858";    D="RUBEN DARIO CORONEL MORALES";    E="2309"; F=46400; G=1000000},
    @{Row=35; C="1065003738"; D="JUAN MANUEL DEGIOVANNI PRECIADO";E="2309"; F=46400; G=1000000},
    @{Row=36; C="60405011";   D="ROSARIO MARIA FLOREZ TEHERAN";   E="2310"; F=46400; G=1000000},
    @{Row=37; C="1047455587"; D="FERNANDO ANDRES OVALLE CORDERO"; E="2310"; F=46400; G=1000000},
    @{Row=38; C="7920858";    D="RUBEN DARIO CORONEL MORALES";    E="2310"; F=46400; G=1000000},
    @{Row=39; C="1065003738"; D="JUAN MANUEL DEGIOVANNI PRECIADO";E="2310"; F=46400; G=1000000},
    @{Row=40; C="60405011";   D="ROSARIO MARIA FLOREZ TEHERAN";   E="2311"; F=37333; G=1000000},
    @{Row=41; C="1047455587"; D="FERNANDO ANDRES OVALLE CORDERO"; E="2311"; F=37333; G=1000000},
    @{Row=42; C="7920858";    D="RUBEN DARIO CORONEL MORALES";    E="2311"; F=37333; G=1000000},
    @{Row=43; C="1065003738"; D="JUAN MANUEL DEGIOVANNI PRECIADO";E="2311"; F=37333; G=1000000}
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
}
